$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("closed bugs in last iteration")

# Update the wording of the "Loading mask..." bug title (row 5, bug 2929769)
$ws.Range("B5").Value = "Loading mask glitch when deleting more packages in Packages grid from Feed details"

# Update the wording of the Grafana bug title (row 12, bug 2901954)
$ws.Range("B12").Value = "Security vulnerabilities in ni-grafana"

# Add a new closed-bug row (row 17, bug 2935004) reusing the formatting of the row above it
$ws.Range("A16:B16").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A17").Value = 2935004
$ws.Range("B17").Value = "Scripts (JupyterHub) not shown in navigation tree when a CNI with NetworkPolicy support is installed"
$ws.Range("C17").Value = "Closed"

# Move the active selection to where it ended up after the edit
$ws.Range("B20").Select()
